$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters in the order values are listed per row (column H / G left untouched)
$cols = @("B","C","D","E","F","I","J","K","L","M","N")

# New vm_pu values per data row (row 2..25), one array per row aligned with $cols
$data = @{
    2 = @(1.02, 1.064724957604878, 1.065175606566462, 1.078354936185895, 1.085505988433052, 1.057525501861298, 1.069683727384349, 1.067889554220751, 1.081033827295637, 1.088166242620205, 1.026653370081924)
    3 = @(1.02, 1.065974041395738, 1.066174040841463, 1.079594950794325, 1.086817963724008, 1.057991001003911, 1.070586716200008, 1.068702878649338, 1.082090679619974, 1.089296237341153, 1.026979799345733)
    4 = @(1.02, 1.066781609236174, 1.066819457602586, 1.080397013076276, 1.087666741834074, 1.058290595151519, 1.071169798661195, 1.069227902708966, 1.082773654718781, 1.090026704884704, 1.027190039827193)
    5 = @(1.02, 1.067120952719469, 1.067090640115097, 1.080734128394972, 1.088023533175576, 1.058416159137786, 1.071414638541483, 1.069448325070317, 1.08306056935761, 1.090333624849693, 1.027278190688714)
    6 = @(1.02, 1.067177920804886, 1.067136164067667, 1.080790727412555, 1.088083438040952, 1.058437219319434, 1.071455731407838, 1.06948531754636, 1.083108731453457, 1.09038514824831, 1.027292977888953)
    7 = @(1.02, 1.066786144177561, 1.066823081746468, 1.08040151790604, 1.087671509431418, 1.058292274456457, 1.071173071352634, 1.069230849171075, 1.082777489300956, 1.090030806623635, 1.027191218622908)
    8 = @(1.02, 1.065147232488627, 1.065513165080281, 1.078774070683184, 1.085949410927949, 1.05768315459373, 1.069989148830758, 1.06816468117633, 1.081391178909235, 1.08854827908806, 1.026763891901978)
    9 = @(1.02, 1.062253970031023, 1.063199964968305, 1.075903810372675, 1.082913522070762, 1.0565973841938, 1.067893551991011, 1.066276293322738, 1.078941485143874, 1.085930291476358, 1.026003345103504)
    10 = @(1.02, 1.060321377631235, 1.06165438130049, 1.073988465432697, 1.080888525508747, 1.055865106337555, 1.066490051462605, 1.065010756097756, 1.077303609186457, 1.084181053335416, 1.025491201769823)
    11 = @(1.02, 1.059483610458116, 1.060984281673373, 1.073158626844331, 1.080011385929148, 1.055546004437451, 1.065880763443341, 1.064461170106137, 1.076593231038147, 1.08342264991068, 1.025268216092983)
    12 = @(1.02, 1.059172281062202, 1.060735246501594, 1.072850312200104, 1.079685528238701, 1.055427170585454, 1.065654209220561, 1.064256786674947, 1.076329186615683, 1.083140795931205, 1.025185204495488)
    13 = @(1.02, 1.059239068914631, 1.060788671327479, 1.072916450213915, 1.079755428093257, 1.055452674686757, 1.065702816662267, 1.06430063863614, 1.076385833159767, 1.083201261368507, 1.02520301914088)
    14 = @(1.02, 1.059457878860899, 1.060963699020945, 1.07313314302185, 1.079984451450701, 1.055536187829415, 1.065862041261658, 1.064444280673003, 1.076571408699537, 1.083399354834481, 1.025261358103063)
    15 = @(1.02, 1.059592675512463, 1.061071522017737, 1.073266644490383, 1.080125553702051, 1.055587602543124, 1.06596011323737, 1.064532751039359, 1.076685724191133, 1.083521386893624, 1.025297278150988)
    16 = @(1.02, 1.060376957305694, 1.061698835458963, 1.074043528659682, 1.080946731696236, 1.055886241368062, 1.066530454740652, 1.065047196400003, 1.077350729799303, 1.084231365327584, 1.025505974728147)
    17 = @(1.02, 1.060868661660052, 1.062092102908892, 1.074530716319153, 1.081461751433956, 1.056073027469316, 1.06688779452925, 1.065369464363678, 1.077767555671831, 1.08467645362862, 1.02563655620766)
    18 = @(1.02, 1.061155373872302, 1.062321407144303, 1.074814838395908, 1.081762124766265, 1.056181781777433, 1.067096074156139, 1.065557283594914, 1.078010570752731, 1.084935972559544, 1.025712604113335)
    19 = @(1.02, 1.061253120072868, 1.062399580123027, 1.074911708937099, 1.081864539502661, 1.056218831181795, 1.06716706665875, 1.065621299011013, 1.078093413641092, 1.085024446037292, 1.025738514484393)
    20 = @(1.02, 1.060815915838614, 1.062049917534727, 1.074478450502126, 1.081406497713185, 1.056053007256363, 1.066849470952262, 1.0653349040306, 1.077722845849343, 1.084628709566193, 1.025622558269068)
    21 = @(1.02, 1.059393448829263, 1.060912161344668, 1.073069334505474, 1.079917011161039, 1.055511603728277, 1.065815160164234, 1.064401988392518, 1.076516766275969, 1.083341025389809, 1.025244183840764)
    22 = @(1.02, 1.05849824440043, 1.060196053612988, 1.072182928059655, 1.078980226621184, 1.055169435261312, 1.065163472337886, 1.063814021814985, 1.075757422809289, 1.082530542548295, 1.025005215261306)
    23 = @(1.02, 1.058972890400184, 1.060575748257112, 1.072652871630149, 1.079476862045121, 1.055350993227874, 1.065509075677604, 1.064125848030119, 1.076160064108966, 1.08296027786647, 1.025131998696955)
    24 = @(1.02, 1.060839749684773, 1.062068979542859, 1.074502067317897, 1.081431464580656, 1.056062054130894, 1.066866788192425, 1.065350520833951, 1.077743048640833, 1.084650283341844, 1.025628883699782)
    25 = @(1.02, 1.063002594596308, 1.063798582258647, 1.076646153400474, 1.083698547913156, 1.056879562131191, 1.068436438584595, 1.066765643818439, 1.079575615364685, 1.086607782098281, 1.026200862719115)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}
